# Insert a new data row at row 578, shifting existing rows 578-629 down to 579-630.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(578).Insert()

$ws.Range("A578").Value = 3
$ws.Range("B578").Value = "Femacal de La Calera"
$ws.Range("C578").Value = "Coquimbo"
$ws.Range("D578").Value = 45166
$ws.Range("D578").NumberFormat = $ws.Range("D579").NumberFormat
$ws.Range("E578").Value = 5
$ws.Range("F578").Value = 100112043
$ws.Range("G578").Value = "Pepino ensalada"
$ws.Range("H578").Value = "Sin especificar"
$ws.Range("I578").Value = "Primera"
$ws.Range("J578").Value = 60
$ws.Range("K578").Value = 7500
$ws.Range("L578").Value = 7500
$ws.Range("M578").Value = 7500
$ws.Range("N578").Value = "$/caja 60 unidades"
$ws.Range("O578").Value = "Región de Arica y Parinacota"
$ws.Range("P578").Value = 125
$ws.Range("Q578").Value = 60
$ws.Range("R578").Value = "Hortaliza"
